$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 28, pushing the existing rows 28-40 down to 30-42.
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()

# --- New row 28 ---
$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Macroferia Regional de Talca"
$ws.Range("C28").Value = "Maule"
$ws.Range("D28").Value = 44489
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100107
$ws.Range("H28").Value = "Otros"
$ws.Range("I28").Value = 100107002
$ws.Range("J28").Value = "Chirimoya"
$ws.Range("K28").Value = "Cultivar IV Región"
$ws.Range("L28").Value = "Especial"
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 27000
$ws.Range("O28").Value = 27000
$ws.Range("P28").Value = 27000
$ws.Range("Q28").Value = "$/bandeja 10 kilos"
$ws.Range("R28").Value = "Provincia de Limarí"
$ws.Range("S28").Value = 2700
$ws.Range("T28").Value = 10

# --- New row 29 ---
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").Value = 44489
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100107
$ws.Range("H29").Value = "Otros"
$ws.Range("I29").Value = 100107002
$ws.Range("J29").Value = "Chirimoya"
$ws.Range("K29").Value = "Cultivar IV Región"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 25000
$ws.Range("O29").Value = 25000
$ws.Range("P29").Value = 25000
$ws.Range("Q29").Value = "$/bandeja 10 kilos"
$ws.Range("R29").Value = "Provincia de Limarí"
$ws.Range("S29").Value = 2500
$ws.Range("T29").Value = 10
